$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K header: "Имя детали на раскладке" (part name shown on the layout),
# added as a shared string and referenced from the new header cell K1.
$ws.Range("K1").Value = "Имя детали на раскладке"

# Column K gets the same (wide) column width as column J, which holds
# "Материал" and is already set to ~72.93 characters wide.
$ws.Columns.Item(11).ColumnWidth = 72.1667

# The edit finished with the new cell selected / scrolled into view.
$ws.Range("K1").Select()
